$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6

# Row 4 updates
$ws.Range("G4").Value = 3.2
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 3.65
$ws.Range("L4").Value = 2.62
$ws.Range("N4").Value = 11.5
$ws.Range("P4").Value = 3.45
$ws.Range("U4").Value = 1.55
$ws.Range("V4").Value = 2.15
$ws.Range("W4").Value = 11.25
$ws.Range("X4").Value = 18.5
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 45
$ws.Range("AA4").Value = 26
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 11.75
$ws.Range("AD4").Value = 6.6
$ws.Range("AJ4").Value = 8.5
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 15.5
$ws.Range("AM4").Value = 22
$ws.Range("AN4").Value = 5.2
$ws.Range("AO4").Value = 17.5
$ws.Range("AP4").Value = 22
$ws.Range("AQ4").Value = 80
$ws.Range("AR4").Value = 100
$ws.Range("AT4").Value = 2.9
$ws.Range("AU4").Value = 6.5
$ws.Range("AW4").Value = 4.15
$ws.Range("AX4").Value = 10.5
$ws.Range("AY4").Value = 16.5
$ws.Range("AZ4").Value = 37
$ws.Range("BA4").Value = 60
$ws.Range("BB4").Value = 175
